$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.579067386262
$ws.Range("C2").Value = 10.18873005904503
$ws.Range("D2").Value = 4.691249409614826
$ws.Range("F2").Value = 26.18503826315785
$ws.Range("G2").Value = 32.59626513152126
$ws.Range("H2").Value = 14.79338156589239
$ws.Range("I2").Value = 22.52923164887649
$ws.Range("L2").Value = 10.67052432394813
$ws.Range("N2").Value = 17.20063336986486
$ws.Range("B3").Value = 16.01564520329605
$ws.Range("C3").Value = 9.842143837713504
$ws.Range("D3").Value = 4.691937392981539
$ws.Range("F3").Value = 26.05486253810532
$ws.Range("G3").Value = 32.30085534803344
$ws.Range("H3").Value = 14.81397825927603
$ws.Range("I3").Value = 22.5933391259103
$ws.Range("L3").Value = 10.64414220510856
$ws.Range("N3").Value = 17.26545740299296
$ws.Range("B4").Value = 15.66260863826682
$ws.Range("C4").Value = 9.621098179482944
$ws.Range("D4").Value = 4.692517565558451
$ws.Range("F4").Value = 25.98378304150948
$ws.Range("G4").Value = 32.13187000479413
$ws.Range("H4").Value = 14.83024514762112
$ws.Range("I4").Value = 22.63879947698332
$ws.Range("L4").Value = 10.63017164491799
$ws.Range("N4").Value = 17.30713762281366
$ws.Range("B5").Value = 15.51719350816915
$ws.Range("C5").Value = 9.529027829678416
$ws.Range("D5").Value = 4.692793921690373
$ws.Range("F5").Value = 25.95706519014462
$ws.Range("G5").Value = 32.06620014351142
$ws.Range("H5").Value = 14.83778195259861
$ws.Range("I5").Value = 22.65885150908445
$ws.Range("L5").Value = 10.62504262546203
$ws.Range("N5").Value = 17.32459636356965
$ws.Range("B6").Value = 15.49296104485675
$ws.Range("C6").Value = 9.513621777747391
$ws.Range("D6").Value = 4.692842229925111
$ws.Range("F6").Value = 25.95276504134843
$ws.Range("G6").Value = 32.05549061618748
$ws.Range("H6").Value = 14.83908819046463
$ws.Range("I6").Value = 22.66227312825355
$ws.Range("L6").Value = 10.62422512649681
$ws.Range("N6").Value = 17.32752402872409
$ws.Range("B7").Value = 15.6606534710065
$ws.Range("C7").Value = 9.619864442249963
$ws.Range("D7").Value = 4.692521130578217
$ws.Range("F7").Value = 25.98341358800148
$ws.Range("G7").Value = 32.1309713357448
$ws.Range("H7").Value = 14.83034311904924
$ws.Range("I7").Value = 22.6390637340317
$ws.Range("L7").Value = 10.63010018452916
$ws.Range("N7").Value = 17.3073711575149
$ws.Range("B8").Value = 16.38640570702853
$ws.Range("C8").Value = 10.07098252542673
$ws.Range("D8").Value = 4.691454042899666
$ws.Range("F8").Value = 26.13833176893217
$ws.Range("G8").Value = 32.4918825247119
$ws.Range("H8").Value = 14.7997304938801
$ws.Range("I8").Value = 22.55006622525195
$ws.Range("L8").Value = 10.66096775437101
$ws.Range("N8").Value = 17.22259580626435
$ws.Range("B9").Value = 17.74406074786492
$ws.Range("C9").Value = 10.88717231303919
$ws.Range("D9").Value = 4.690601853364161
$ws.Range("F9").Value = 26.51116305461703
$ws.Range("G9").Value = 33.29429532275172
$ws.Range("H9").Value = 14.76852300151666
$ws.Range("I9").Value = 22.42422048515957
$ws.Range("L9").Value = 10.73898834357076
$ws.Range("N9").Value = 17.07118482948501
$ws.Range("B10").Value = 18.6905645565574
$ws.Range("C10").Value = 11.44154896410685
$ws.Range("D10").Value = 4.690717183824213
$ws.Range("F10").Value = 26.82533877005786
$ws.Range("G10").Value = 33.93610647365009
$ws.Range("H10").Value = 14.76327967339237
$ws.Range("I10").Value = 22.36181589656046
$ws.Range("L10").Value = 10.80668472958476
$ws.Range("N10").Value = 16.96889006112207
$ws.Range("B11").Value = 19.10814363129902
$ws.Range("C11").Value = 11.68331209941453
$ws.Range("D11").Value = 4.690927580117821
$ws.Range("F11").Value = 26.97654495538608
$ws.Range("G11").Value = 34.23816607066667
$ws.Range("H11").Value = 14.76475100661156
$ws.Range("I11").Value = 22.34002548104129
$ws.Range("L11").Value = 10.83966254058994
$ws.Range("N11").Value = 16.92427593401417
$ws.Range("B12").Value = 19.26426409299292
$ws.Range("C12").Value = 11.7733201068988
$ws.Range("D12").Value = 4.69102970481325
$ws.Range("F12").Value = 27.03495092990875
$ws.Range("G12").Value = 34.35388819800219
$ws.Range("H12").Value = 14.7658633851202
$ws.Range("I12").Value = 22.33272804118275
$ws.Range("L12").Value = 10.85245734662993
$ws.Range("N12").Value = 16.90765636064327
$ws.Range("B13").Value = 19.23073209924918
$ws.Range("C13").Value = 11.7540045009748
$ws.Range("D13").Value = 4.691006715949757
$ws.Range("F13").Value = 27.0223218535539
$ws.Range("G13").Value = 34.32890765554713
$ws.Range("H13").Value = 14.76559911515137
$ws.Range("I13").Value = 22.33425716030084
$ws.Range("L13").Value = 10.84968822642921
$ws.Range("N13").Value = 16.91122348064615
$ws.Range("B14").Value = 19.12102873387341
$ws.Range("C14").Value = 11.69074825274226
$ws.Range("D14").Value = 4.690935533474307
$ws.Range("F14").Value = 26.98132731472926
$ws.Range("G14").Value = 34.24766042692101
$ws.Range("H14").Value = 14.76483139304711
$ws.Range("I14").Value = 22.33940596525801
$ws.Range("L14").Value = 10.84070906948055
$ws.Range("N14").Value = 16.9229031311529
$ws.Range("B15").Value = 19.05356693342136
$ws.Range("C15").Value = 11.65179991630386
$ws.Range("D15").Value = 4.690894848462826
$ws.Range("F15").Value = 26.95636505804449
$ws.Range("G15").Value = 34.19806509183933
$ws.Range("H15").Value = 14.76443345888852
$ws.Range("I15").Value = 22.34268416255673
$ws.Range("L15").Value = 10.83524881863503
$ws.Range("N15").Value = 16.93009300666969
$ws.Range("B16").Value = 18.66300106279264
$ws.Range("C16").Value = 11.42553555500378
$ws.Range("D16").Value = 4.690706589798068
$ws.Range("F16").Value = 26.81562012329722
$ws.Range("G16").Value = 33.91655972485922
$ws.Range("H16").Value = 14.76326118686973
$ws.Range("I16").Value = 22.3633731846789
$ws.Range("L16").Value = 10.80457287949266
$ws.Range("N16").Value = 16.97184422813099
$ws.Range("B17").Value = 18.4199697077158
$ws.Range("C17").Value = 11.2840268184839
$ws.Range("D17").Value = 4.690631364922866
$ws.Range("F17").Value = 26.73136945566542
$ws.Range("G17").Value = 33.74637163897638
$ws.Range("H17").Value = 14.76353038359522
$ws.Range("I17").Value = 22.37775903907767
$ws.Range("L17").Value = 10.78630831366932
$ws.Range("N17").Value = 16.99794814454604
$ws.Range("B18").Value = 18.27896904236733
$ws.Range("C18").Value = 11.20165562840855
$ws.Range("D18").Value = 4.690602978107441
$ws.Range("F18").Value = 26.68369344462
$ws.Range("G18").Value = 33.64944126100711
$ws.Range("H18").Value = 14.76404820209205
$ws.Range("I18").Value = 22.38665412066067
$ws.Range("L18").Value = 10.77600888413874
$ws.Range("N18").Value = 17.01314324907522
$ws.Range("B19").Value = 18.23102441461011
$ws.Range("C19").Value = 11.17359940960717
$ws.Range("D19").Value = 4.690595930487219
$ws.Range("F19").Value = 26.66768690125643
$ws.Range("G19").Value = 33.61679012160892
$ws.Range("H19").Value = 14.76428584284612
$ws.Range("I19").Value = 22.38977228027185
$ws.Range("L19").Value = 10.77255723754317
$ws.Range("N19").Value = 17.01831914687291
$ws.Range("B20").Value = 18.44596767028886
$ws.Range("C20").Value = 11.29919233069961
$ws.Range("D20").Value = 4.690637834522581
$ws.Range("F20").Value = 26.7402573824247
$ws.Range("G20").Value = 33.76439018465602
$ws.Range("H20").Value = 14.7634641536421
$ws.Range("I20").Value = 22.37616335782589
$ws.Range("L20").Value = 10.78823135079677
$ws.Range("N20").Value = 16.99515063314381
$ws.Range("B21").Value = 19.15330679064618
$ws.Range("C21").Value = 11.70937032415984
$ws.Range("D21").Value = 4.690955834192908
$ws.Range("F21").Value = 26.99333762660835
$ws.Range("G21").Value = 34.2714892892786
$ws.Range("H21").Value = 14.76504182009383
$ws.Range("I21").Value = 22.33786770149807
$ws.Range("L21").Value = 10.84333819645577
$ws.Range("N21").Value = 16.91946508622625
$ws.Range("B22").Value = 19.60383813347481
$ws.Range("C22").Value = 11.96843865739739
$ws.Range("D22").Value = 4.69129441566892
$ws.Range("F22").Value = 27.16540631893333
$ws.Range("G22").Value = 34.61065687839897
$ws.Range("H22").Value = 14.76930915128046
$ws.Range("I22").Value = 22.31840275742044
$ws.Range("L22").Value = 10.88113855054364
$ws.Range("N22").Value = 16.87160152333675
$ws.Range("B23").Value = 19.36449866477936
$ws.Range("C23").Value = 11.83100596563115
$ws.Range("D23").Value = 4.691101829454251
$ws.Range("F23").Value = 27.07297529107569
$ws.Range("G23").Value = 34.42896569793659
$ws.Range("H23").Value = 14.76673537327058
$ws.Range("I23").Value = 22.32828085445833
$ws.Range("L23").Value = 10.86080290626393
$ws.Range("N23").Value = 16.8970011134452
$ws.Range("B24").Value = 18.43421796347139
$ws.Range("C24").Value = 11.29233916380678
$ws.Range("D24").Value = 4.690634863279232
$ws.Range("F24").Value = 26.73623677743236
$ws.Range("G24").Value = 33.75624114750281
$ws.Range("H24").Value = 14.76349296530075
$ws.Range("I24").Value = 22.37688282047963
$ws.Range("L24").Value = 10.78736131874647
$ws.Range("N24").Value = 16.99641480349179
$ws.Range("B25").Value = 17.38500817758171
$ws.Range("C25").Value = 10.67409474020449
$ws.Range("D25").Value = 4.690701133481418
$ws.Range("F25").Value = 26.4030863295341
$ws.Range("G25").Value = 33.06762456633503
$ws.Range("H25").Value = 14.77386652094835
$ws.Range("I25").Value = 22.45301311293609
$ws.Range("L25").Value = 10.71603694442111
$ws.Range("N25").Value = 17.11056748786224
